$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-4 are cyclically rotated:
#   new row2 <- old row3
#   new row3 <- old row4
#   new row4 <- old row2
# Only columns A, B, D, E, F, G, H, Q, R actually change.

$ws.Range("A2").Value = 111463862
$ws.Range("B2").Value = 89369
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 5447
$ws.Range("F2").Value = "Vedticka"
$ws.Range("G2").Value = "Fuscoporia viticola"
$ws.Range("H2").Value = "(Schwein.) Murrill"
$ws.Range("Q2").Value = 554109.1038748255
$ws.Range("R2").Value = 7007938.027731327

$ws.Range("A3").Value = 111463670
$ws.Range("B3").Value = 96674
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 219880
$ws.Range("F3").Value = "Kransrams"
$ws.Range("G3").Value = "Polygonatum verticillatum"
$ws.Range("H3").Value = "(L.) All."
$ws.Range("Q3").Value = 554151.0634843309
$ws.Range("R3").Value = 7007942.793868498

$ws.Range("A4").Value = 111463857
$ws.Range("B4").Value = 89405
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = "Ullticka"
$ws.Range("G4").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H4").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q4").Value = 554109.1038748255
$ws.Range("R4").Value = 7007938.027731327
